$p = $ppt.ActivePresentation

# --- Slide 5: "Timeline & Milestones" table - adjust column widths ---
$s1 = $p.Slides.Item(5)
$tbl1 = $s1.Shapes.Item(3).Table
$tbl1.Columns.Item(1).Width = 871093 / 12700
$tbl1.Columns.Item(3).Width = 1306639 / 12700
$tbl1.Columns.Item(4).Width = 4355466 / 12700

# --- Slide 8: "Investment Summary" table - adjust column widths ---
$s2 = $p.Slides.Item(8)
$tbl2 = $s2.Shapes.Item(3).Table
$tbl2.Columns.Item(1).Width = 1742186 / 12700
$tbl2.Columns.Item(2).Width = 1045311 / 12700
$tbl2.Columns.Item(3).Width = 2003514 / 12700
$tbl2.Columns.Item(4).Width = 1132421 / 12700
$tbl2.Columns.Item(5).Width = 871093 / 12700
$tbl2.Columns.Item(6).Width = 871093 / 12700
$tbl2.Columns.Item(7).Width = 1045311 / 12700
